# "improved excel result sheet"
# - header row (A1:D1) gets a green solid fill + centered alignment
# - active selection moves back to A1
# - header row is made taller (custom height)
# - data columns A-D are narrowed slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row formatting (A1:D1) ---------------------------------------
$header = $ws.Range("A1:D1")

# Solid fill: foreground FF66CC00 (green), background FF339966 (teal-green)
$header.Interior.Pattern = 1
$header.Interior.Color = 52326
$header.Interior.PatternColor = 6723891

# Centre the header text both horizontally and vertically
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# --- Row height -------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 16.85

# --- Column widths (characters) ---------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 42
$ws.Columns.Item(2).ColumnWidth = 13.3333333333333
$ws.Columns.Item(3).ColumnWidth = 14.1666666666667
$ws.Columns.Item(4).ColumnWidth = 14

# --- Selection ---------------------------------------------------------------
[void]$ws.Range("A1").Select()
